$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2027-08-12"
$ws.Range("B2").Style = "Normal"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2027-08-12"
$ws.Range("B3").Style = "Normal"
$ws.Range("E5").Value = "❌ EXPIRED 3776 days ago"
$ws.Range("E9").Value = "⚠️ Expires in 22 days"
$ws.Range("E10").Value = "⚠️ Expires in 30 days"
